$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing the existing rows 11-29 down to 12-30
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly record
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44935
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112044
$ws.Range("G11").Value = "Perejil"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 78
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 3000
$ws.Range("N11").Value = "$/docena de atados (3 kilos)"
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 1000
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
